# Weekly price-data refresh: insert one new record as row 85 in the
# "Poroto verde" sheet, pushing the existing rows 85-129 down to 86-130.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 85; this shifts rows
# 85..129 down to 86..130 (carrying all their data/styles with them),
# exactly matching the target diff.
$ws.Rows("85").Insert()

# Populate the freshly-inserted row 85 with the new weekly record.
$ws.Range("A85").Value = 4
$ws.Range("B85").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C85").Value = "Los Lagos"
$ws.Range("D85").Value = 44992
$ws.Range("E85").Value = 10
$ws.Range("F85").Value = 100112031
$ws.Range("G85").Value = "Poroto verde"
$ws.Range("H85").Value = "Magnum"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 45
$ws.Range("K85").Value = 32000
$ws.Range("L85").Value = 32000
$ws.Range("M85").Value = 32000
$ws.Range("N85").Value = "$/saco 25 kilos"
$ws.Range("O85").Value = "Región Metropolitana"
$ws.Range("P85").Value = 1280
$ws.Range("Q85").Value = 25
$ws.Range("R85").Value = "Hortaliza"
